$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 403 (shifts existing rows 403:496 down to 404:497,
# Excel also pushes the dimension ref out to R497 and carries the D-column
# date style down onto the freshly inserted row).
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row 403 with the new data point. The row
# mirrors the existing "Feria Lagunitas de Puerto Montt" / "Pepino ensalada"
# records, only the date, volume, prices and $/Kg differ.
$ws.Cells.Item(403, 1).Value = 4
$ws.Cells.Item(403, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(403, 3).Value = "Los Lagos"
$ws.Cells.Item(403, 4).Value = 45173
$ws.Cells.Item(403, 5).Value = 10
$ws.Cells.Item(403, 6).Value = 100112043
$ws.Cells.Item(403, 7).Value = "Pepino ensalada"
$ws.Cells.Item(403, 8).Value = "Sin especificar"
$ws.Cells.Item(403, 9).Value = "Primera"
$ws.Cells.Item(403, 10).Value = 120
$ws.Cells.Item(403, 11).Value = 15000
$ws.Cells.Item(403, 12).Value = 15000
$ws.Cells.Item(403, 13).Value = 15000
$ws.Cells.Item(403, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(403, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(403, 16).Value = 250
$ws.Cells.Item(403, 17).Value = 60
$ws.Cells.Item(403, 18).Value = "Hortaliza"
